$wb = $excel.ActiveWorkbook

# --- Sheet "general" ---
$ws = $wb.Worksheets.Item("general")
$ws.Range("B3").Value = 94.4371278619396
$ws.Range("B4").Value = 0.01199984550476074
$ws.Range("B6").Value = 33.0671278619396
$ws.Range("B10").Value = 61.37

# --- Sheet "x" ---
$ws = $wb.Worksheets.Item("x")
$ws.Range("B2").Value = 1
$ws.Range("B3").Value = 3
$ws.Range("B4").Value = 5
$ws.Range("B6").Value = 12
$ws.Range("B7").Value = 6
$ws.Range("B10").Value = 2

# --- Sheet "U" ---
$ws = $wb.Worksheets.Item("U")
$ws.Range("B3").Value = 3

# --- Sheet "TBar" ---
$ws = $wb.Worksheets.Item("TBar")
$ws.Range("B4").Value = 24.69770569366316
$ws.Range("B5").Value = 20
$ws.Range("B8").Value = 20.34885527085025
$ws.Range("B14").Value = 27.27819014430416

# --- Sheet "Q" ---
$ws = $wb.Worksheets.Item("Q")
$ws.Range("C7").Value = 109.9450000000008
$ws.Range("C8").Value = 117.5900000000008
$ws.Range("C9").Value = 113.2700000000008
$ws.Range("C10").Value = 119.1550000000008
$ws.Range("C11").Value = 115.8050000000008
$ws.Range("C12").Value = 235.775
$ws.Range("C13").Value = 229.025
$ws.Range("C14").Value = 213.42
$ws.Range("C15").Value = 226.76
$ws.Range("C16").Value = 221.56
$ws.Range("C17").Value = 46.91999999999942
$ws.Range("C18").Value = 36.10499999999942
$ws.Range("C19").Value = 34.91499999999942
$ws.Range("C20").Value = 37.48999999999942
$ws.Range("C21").Value = 39.43499999999941
$ws.Range("C32").Value = 154.3
$ws.Range("C33").Value = 148.3449999999993
$ws.Range("C34").Value = 128.7049999999993
$ws.Range("C35").Value = 146.3249999999993
$ws.Range("C36").Value = 134.2149999999993
$ws.Range("C37").Value = 193.0200000000017
$ws.Range("C38").Value = 202.3100000000017
$ws.Range("C39").Value = 191.2450000000017
$ws.Range("C40").Value = 208.9250000000017
$ws.Range("C41").Value = 197.6600000000017
$ws.Range("C52").Value = 250.970000000001
$ws.Range("C53").Value = 260.9900000000009
$ws.Range("C54").Value = 252.975000000001
$ws.Range("C55").Value = 269.580000000001
$ws.Range("C56").Value = 250.575000000001
$ws.Range("C62").Value = 235.775
$ws.Range("C63").Value = 229.025
$ws.Range("C64").Value = 213.42
$ws.Range("C65").Value = 226.76
$ws.Range("C66").Value = 221.56
$ws.Range("C67").Value = 250.970000000001
$ws.Range("C68").Value = 260.9900000000009
$ws.Range("C69").Value = 252.975000000001
$ws.Range("C70").Value = 269.580000000001
$ws.Range("C71").Value = 250.575000000001

# --- Sheet "L" ---
$ws = $wb.Worksheets.Item("L")
$ws.Range("C7").Value = 0
$ws.Range("C8").Value = 0
$ws.Range("C9").Value = 0
$ws.Range("C10").Value = 0
$ws.Range("C11").Value = 0
